# Update the "Förändrad" date column (C) for rows 2-18 from 45183 to 45184
# (i.e. increment the serial date value by 1 day, 2023-09-14 -> 2023-09-15)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 18; $row++) {
    $cell = $ws.Cells.Item($row, 3)  # Column C
    if ($cell.Value2 -eq 45183) {
        $cell.Value2 = 45184
    }
}
